$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Second iteration of the receptive-field calculation: Block 1's first
# --- Conv2d layer now uses a dilated kernel (dilation 1 -> 2). All the
# --- downstream formulas (J9:K23 etc.) recompute automatically because
# --- they reference F8 through the shared FLOOR(...) formula chain.
$ws.Range("F8").Value = 2

# Scroll the view down a bit and leave the selection on I26, matching the
# cursor position the author ended up on after re-running the sheet.
$ws.Application.Goto($ws.Range("A4"), $false)
$ws.Range("I26").Select()

# Record the newly observed "Max accuracy" for this run, now formatted
# with two decimal places (0.00%) instead of the previous whole-number %.
$ws.Range("M25").Value = 0.8023
$ws.Range("M25").NumberFormat = "0.00%"
